$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Race Unknown" / "Race Other" header columns (CV1 <-> CW1)
# to reorder the race options for consistency.
$ws.Range("CV1").Value = "Race Other"
$ws.Range("CW1").Value = "Race Unknown"
